# Adds new question-variant lines to the "question" column (C) of several
# rows in the NPCEditor_data sheet. Each target cell already holds a
# newline-separated list of question paraphrases; we append additional
# paraphrases to the end of the existing list, preserving everything that
# was already there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Append-Questions {
    param(
        [string]$CellRef,
        [string[]]$NewLines
    )
    $cell = $ws.Range($CellRef)
    $current = $cell.Value()
    $addition = [string]::Join("`n", $NewLines)
    $cell.Value = $current + "`n" + $addition
}

Append-Questions "C13"  @("what do you do in the navy?", "what is your Navy job?")
Append-Questions "C39"  @("Was adjusting to the Navy tough?", "Was the first year hard?")
Append-Questions "C65"  @("do you play video games?")
Append-Questions "C73"  @("what do you like about the Navy?")
Append-Questions "C104" @("tell me a story.")
Append-Questions "C131" @("why did you study computer science?")
Append-Questions "C132" @("was the Navy your first choice?", "How did you decide to join the navy?")
Append-Questions "C147" @("Have you ever swim in the middle of the ocean?", "have you swam in the ocean?", "do you get to swim in the ocean?", "when have you swam somewhere cool?")
Append-Questions "C166" @("What is the first year like in the Navy?", "Was the first year hard?")
Append-Questions "C174" @("Tell me about navy life.", "what is Navy life like?")
Append-Questions "C221" @("who did you idolize growing up?", "did you look up to someone as a kid or student?")
Append-Questions "C228" @("was your job dangerous at all?", "did you do anything risky sometime?")
Append-Questions "C275" @("do you enjoy being in the navy?", "do you like the Navy?")
Append-Questions "C301" @("what will you do after the navy?", "what will you do after you graduate?", "what will you do after college?")
Append-Questions "C308" @("what is the best job out there?")
Append-Questions "C321" @("What should I do?")
Append-Questions "C350" @("did you like college?", "how was college?", "tell me about college.", "i want to know about college.", "do you like college?")
Append-Questions "C361" @("Did you ever experience combat?")

Write-Host "Done appending question variants."
